$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.379.11"
$ws.Range("E2").Value = "  +3.99%  "

$ws.Range("D3").Value = "2.262.29"
$ws.Range("E3").Value = "  +1.86%  "

$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").Value = "'231.16"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("E6").Value = "  +0.76%  "

$ws.Range("D7").Value = "'60.99"
$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").Value = "'0.410"
$ws.Range("E9").Value = "  +2.65%  "

$ws.Range("D10").Value = "'0.0914"
$ws.Range("E10").Value = "  +2.94%  "

$ws.Range("E11").Value = "  +0.76%  "

$ws.Range("D12").Value = "2.604.77"
$ws.Range("E12").Value = "  +2.09%  "

$ws.Range("D13").Value = "'15.64"
$ws.Range("E13").Value = "  +0.03%  "

$ws.Range("D14").Value = "'22.49"
$ws.Range("E14").Value = "  +3.69%  "

$ws.Range("D15").Value = "'5.68"
$ws.Range("E15").Value = "  +2.38%  "

$ws.Range("D16").Value = "'0.805"
$ws.Range("E16").Value = "  +1.14%  "

$ws.Range("D17").Value = "2.274.65"
$ws.Range("E17").Value = "  +2.42%  "

$ws.Range("D18").Value = "43.297.71"
$ws.Range("E18").Value = "  +4.24%  "

$ws.Range("D19").Value = "0.0₃0929"
$ws.Range("E19").Value = "  +3.93%  "

$ws.Range("D20").Value = "'72.89"
$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("D21").Value = "'6.20"
$ws.Range("E21").Value = "  +2.73%  "

$ws.Range("D22").Value = "'247.89"
$ws.Range("E22").Value = "  -0.82%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'2.57"
$ws.Range("E24").Value = "  +8.07%  "

$ws.Range("D25").Value = "'2.38"
$ws.Range("E25").Value = "  +4.31%  "

$ws.Range("D26").Value = "'9.74"
$ws.Range("E26").Value = "  +1.95%  "

$ws.Range("D27").Value = "'169.63"
$ws.Range("E27").Value = "  +1.17%  "

$ws.Range("D28").Value = "'0.142"
$ws.Range("E28").Value = "  +1.30%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'20.53"
$ws.Range("E29").Value = "  +3.01%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'1.48"
$ws.Range("E30").Value = "  +4.75%  "

$ws.Range("D31").Value = "'2.67"
$ws.Range("E31").Value = "  +1.48%  "

$ws.Range("E32").Value = "  -0.89%  "

$ws.Range("D33").Value = "'5.01"
$ws.Range("E33").Value = "  +1.86%  "

$ws.Range("D34").Value = "'4.71"
$ws.Range("E34").Value = "  +2.47%  "

$ws.Range("D35").Value = "'0.0652"
$ws.Range("E35").Value = "  +5.03%  "

$ws.Range("D36").Value = "'6.42"
$ws.Range("E36").Value = "  -2.20%  "

$ws.Range("D37").Value = "'2.39"
$ws.Range("E37").Value = "  +1.56%  "

$ws.Range("D38").Value = "'3.58"
$ws.Range("E38").Value = "  -2.61%  "

$ws.Range("D39").Value = "'0.0250"
$ws.Range("E39").Value = "  +4.96%  "

$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("D41").Value = "'8.65"
$ws.Range("E41").Value = "  +0.76%  "

$ws.Range("D42").Value = "'0.000219"
$ws.Range("E42").Value = "  -10.47%  "

$ws.Range("D43").Value = "'0.0969"
$ws.Range("E43").Value = "  -1.06%  "

$ws.Range("D44").Value = "'1.20"
$ws.Range("E44").Value = "  -0.40%  "

$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'4.38"
$ws.Range("E45").Value = "  -9.84%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'96.74"
$ws.Range("E46").Value = "  -1.85%  "

$ws.Range("D47").Value = "1.462.66"
$ws.Range("E47").Value = "  -0.08%  "

$ws.Range("D48").Value = "'16.70"
$ws.Range("E48").Value = "  +1.43%  "

$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "'2.76"
$ws.Range("E49").Value = "  -1.71%  "

$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "'1.08"
$ws.Range("E50").Value = "  +0.39%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.478.23"
$ws.Range("E51").Value = "  +1.83%  "
